$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.717000000000001
$ws.Range("E4").Value = 13.265

$ws.Range("E5").Value = 13.254

$ws.Range("D7").Value = -7.612

$ws.Range("E8").Value = 13.718

$ws.Range("D16").Value = -8.219999999999999
$ws.Range("E16").Value = 12.978
